# Regenerate the "K" column (column G) values for rows 2-37.
# The previous values were derived from an old "Strike#" calculation;
# this writes the recalculated strikeout (K) values (s_vals) in their place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(4,6,4,5,5,6,7,5,8,6,5,5,0,10,1,8,3,1,4,2,5,5,6,4,7,3,1,6,7,6,2,6,2,3,3,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
